$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.001.57"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.909.64"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "0.7924"
$ws.Range("E5").Value = "  +5.78%  "
$ws.Range("D6").Value = "242.00"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +2.83%  "
$ws.Range("D9").Value = "26.32"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("D10").Value = "0.06892"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "0.07996"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "1.904.72"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "0.7438"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").Value = "5.188"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "93.06"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "29.998.57"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "5.868"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("D19").Value = "245.83"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").Value = "0.000007741"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "2.151.13"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "6.836"
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("D25").Value = "168.08"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "9.233"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "0.1401"
$ws.Range("E27").Value = "  +9.87%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "2.031"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").Value = "1.516"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "4.315"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "0.05550"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "4.079"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "0.7339"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").Value = "2.721"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "0.01924"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "2.783"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").Value = "6.146"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").Value = "0.4417"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "72.32"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "0.8370"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Value = "1.876"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").Value = "100.53"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "7.547"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "987.69"
$ws.Range("E48").Value = "  +8.27%  "
$ws.Range("D49").Value = "2.054.33"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "36.25"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "1.478"
$ws.Range("E51").Value = "  -0.13%  "
